$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp in the title cell (row 1)
$ws.Range("A1").Value = "Datos actualizados a 11 de Abril de 2020 a las 22:22"

# Helper to write a full data row: country name in column A, then B..H values
function Set-CountryRow {
    param(
        [int]$Row,
        [string]$Name,
        $Values
    )
    $ws.Cells.Item($Row, 1).Value = $Name
    for ($i = 0; $i -lt $Values.Length; $i++) {
        $ws.Cells.Item($Row, 2 + $i).Value = $Values[$i]
    }
}

# Row 4: Estados Unidos - updated totals
Set-CountryRow 4 "Estados Unidos" @(525701, 22825, 28755, 476634, 11053, 1565, 20312)

# Row 8: Alemania - updated totals
Set-CountryRow 8 "Alemania" @(124288, 2117, 57400, 64152, 4895, 0, 2736)

# Row 17: Brasil - updated totals
Set-CountryRow 17 "Brasil" @(20727, 938, 173, 19430, 296, 56, 1124)

# Rows 47/48: Republica Dominicana overtakes Catar
Set-CountryRow 47 "Republica Dominicana" @(2759, 139, 108, 2516, 147, 9, 135)
Set-CountryRow 48 "Catar" @(2728, 216, 247, 2475, 37, 0, 6)

# Rows 57/58: Egipto overtakes Argelia
Set-CountryRow 57 "Egipto" @(1939, 145, 426, 1367, 0, 11, 146)
Set-CountryRow 58 "Argelia" @(1825, 64, 460, 1090, 46, 19, 275)

# Row 73: Armenia - updated totals
Set-CountryRow 73 "Armenia" @(967, 30, 173, 781, 30, 1, 13)

# Row 143: Uganda - updated totals
Set-CountryRow 143 "Uganda" @(53, 0, 4, 49, 0, 0, 0)

# Rows 168/169: Sudan overtakes Angola
Set-CountryRow 168 "Sudan" @(19, 2, 2, 15, 0, 0, 2)
Set-CountryRow 169 "Angola" @(19, 0, 4, 13, 0, 0, 2)
